$d = $word.ActiveDocument

# ============================================================
# Region 3 (bottom): remove cached lastRenderedPageBreak from
# the "This brings us to my final point..." paragraph (it is
# logically relocated earlier, onto the "In current scenario"
# paragraph, by the edit).
# ============================================================
$d.Content.Find.Execute(
    "This brings us to my final point that the given solution to this idea uses AI extraction bot for document text extraction, reinforcement bot for giving competitiveness to gaming participants while doing data keying and vision bots for auto data entry.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This brings us to my final point that the given solution to this idea uses AI extraction bot for document text extraction, reinforcement bot for giving competitiveness to gaming participants while doing data keying and vision bots for auto data entry.",
    2)

# ============================================================
# Region 2 (Slide 3 area)
# ============================================================

# "Another intuition is survey before..." -> "Entertain intuition from the surveys before..."
$d.Content.Find.Execute(
    "Another intuition is survey before in youtube videos where the viewer is asked to fill survey details to proceed watching videos. We can greatly leverage it by keying the data before proceeding to videos/on demand tv shows.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Entertain intuition from the surveys before in youtube videos where the viewer is asked to fill survey details to proceed watching videos. We can greatly leverage it to our idea by keying the data before proceeding to videos/on demand tv shows.",
    2)

# "One of the reference for crowdsorcing platform is amazon mturk..." -> "The crowdsorcing intuition comes from the platform,  amazon mturk..."
$d.Content.Find.Execute(
    "One of the reference for crowdsorcing platform is amazon mturk, where requesters post the task such as question answering, data labeling and so on. So the workers, they search and look into the task whichever suits them and complete within a given timeframe. And in return, they will get respective dollar equivalent of completed task. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The crowdsorcing intuition comes from the platform,  amazon mturk, where requesters post the task such as question answering, data labeling and so on. So the workers, they search and look into the task whichever suits them and complete within a given timeframe. And in return, they will get respective dollar equivalent of completed task. ",
    2)

# "This is not a new idea though, a subset of this idea..." -> "So the idea of 'intralyst sorcing' comes from the intuition , where subset..."
$d.Content.Find.Execute(
    "This is not a new idea though, a subset of this idea has been already implemented to one of our customer named adcb bank in mid-east for women empowerment program. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "So the idea of 'intralyst sorcing' comes from the intuition , where subset of this idea has been already implemented to one of our customer named adcb bank in mid-east for women empowerment program. The process basically randomize the snippet of original bank forms with identification number and those snippets will be given to housewife for data entry so that no will get to know PII details about the bank customers.",
    2)

Write-Host "region 2 done"

# ============================================================
# Region 1 (Slide 2 area)
# ============================================================

# Remove the existing _GoBack bookmark; it is re-added later at
# its new location (end of the "Why I said this..." paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Replace "Not only it motivates participants..." paragraph with
# the new "Eventhough data entry is linked..." paragraph text,
# then append three brand-new paragraphs after it.
$d.Content.Find.Execute(
    "Not only it motivates participants, it also basically demarcates the data with process knowledge. When I say demarcates, the respective process documents can be keyed by any analyst in Sutherland.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Eventhough data entry is linked with process knowledge, it basically demarcates the data with process knowledge. When I say demarcates, the respective process documents can be keyed by any analyst in Sutherland. With the current social media trends, analyst can get the gist about the process knowledge of respective documents.",
    2)

$found = $d.Content.Find.Execute("Eventhough data entry is linked with process knowledge", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$targetPara = $d.Content.Paragraphs.Last
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Eventhough data entry is linked")) {
        $targetPara = $p
    }
}
$targetPara.Range.InsertParagraphAfter()
$targetPara.Next().Range.Text = "The documents can be any usable data entry documents like invoices, scanned PDF, bank forms, insurance , claim forms, healthcare doc and so on."
$targetPara.Next().Range.InsertParagraphAfter()
$targetPara.Next().Next().Range.Text = "Our idea greatly reduces the respective process entry team by assign data keying to any analyst in Sutherland."
$targetPara.Next().Next().Range.InsertParagraphAfter()
$targetPara.Next().Next().Next().Range.Text = "Also, it gives transformation from boring task into a gamified and entertaining way which triggers the internal motive of people participating. "

Write-Host "region1 step A done"

# Replace "Over a period, data keying from documents would be
# monotonous which in turn become a boring and mundane task for
# any human. So the idea gives transformation..." with the new
# "Over a period... L.H.S = R.H.S ... which becomes boring..."
# paragraph text (trailing "So the idea gives transformation..."
# content is relocated to the new "Also, it gives..." paragraph
# created above, so it is dropped here).
$d.Content.Find.Execute(
    "Over a period, data keying from documents would be monotonous which in turn become a boring and mundane task for any human. So the idea gives transformation from boring task into a gamified and entertaining way which triggers the internal motive of people participating. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Over a period, data keying from documents would be a monotonous task where the people select a task and do data emendation and this process continues like L.H.S = R.H.S , which becomes boring for any human. ",
    2)

Write-Host "region1 step B done"

# Insert the two brand-new paragraphs ("We all know..." and "Why
# I said this...") right after the "Let's walk to the problem
# statement." paragraph.
$letsWalkPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Let") -and $p.Range.Text.Contains("walk to the problem statement")) {
        $letsWalkPara = $p
    }
}
$letsWalkPara.Range.InsertParagraphAfter()
$letsWalkPara.Next().Range.Text = "We all know that it would be more fun to spend our waking hours doing recreational things right."
$letsWalkPara.Next().Range.InsertParagraphAfter()
$letsWalkPara.Next().Next().Range.Text = "Why I said this….because the problem statement is related to this…"

Write-Host "region1 step C done"

# Re-add the _GoBack bookmark at the end of the "Why I said
# this..." paragraph (collapsed / zero-length, as it was before).
$whyPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Why I said this")) {
        $whyPara = $p
    }
}
$bmStart = $whyPara.Range.End - 1
$d.Bookmarks.Add("_GoBack", $d.Range($bmStart, $bmStart))

Write-Host "region1 step D done"
